$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update AgTests (F) and AgPosit (G) values for rows 363-409 with refreshed data
$ws.Range("F363").Value = 188334
$ws.Range("G363").Value = 2754
$ws.Range("F364").Value = 168018
$ws.Range("F365").Value = 184578
$ws.Range("G365").Value = 2394
$ws.Range("F366").Value = 338795
$ws.Range("G366").Value = 2841
$ws.Range("F367").Value = 765329
$ws.Range("G367").Value = 3918
$ws.Range("F368").Value = 346892
$ws.Range("G368").Value = 2310
$ws.Range("F370").Value = 180361
$ws.Range("G370").Value = 2043
$ws.Range("F371").Value = 159832
$ws.Range("G371").Value = 1956
$ws.Range("F372").Value = 179434
$ws.Range("G372").Value = 1865
$ws.Range("F373").Value = 349187
$ws.Range("G373").Value = 2368
$ws.Range("F374").Value = 772554
$ws.Range("G374").Value = 3419
$ws.Range("F375").Value = 352034
$ws.Range("G375").Value = 1860
$ws.Range("F377").Value = 176710
$ws.Range("G377").Value = 1826
$ws.Range("F378").Value = 156891
$ws.Range("G378").Value = 1549
$ws.Range("F379").Value = 179958
$ws.Range("G379").Value = 1619
$ws.Range("F380").Value = 344788
$ws.Range("G380").Value = 2024
$ws.Range("F381").Value = 744951
$ws.Range("G381").Value = 2685
$ws.Range("F387").Value = 351483
$ws.Range("F388").Value = 730057
$ws.Range("G388").Value = 2199
$ws.Range("F391").Value = 176185
$ws.Range("F392").Value = 220708
$ws.Range("F393").Value = 306872
$ws.Range("F395").Value = 749213
$ws.Range("F396").Value = 164187
$ws.Range("F397").Value = 108221
$ws.Range("F398").Value = 297730
$ws.Range("F399").Value = 200378
$ws.Range("F400").Value = 149614
$ws.Range("F401").Value = 273034
$ws.Range("F402").Value = 713233
$ws.Range("G402").Value = 1379
$ws.Range("F404").Value = 224576
$ws.Range("G404").Value = 906
$ws.Range("F405").Value = 173283
$ws.Range("F406").Value = 170293
$ws.Range("G406").Value = 674
$ws.Range("F407").Value = 156146
$ws.Range("G407").Value = 668
$ws.Range("F408").Value = 292641
$ws.Range("G408").Value = 861
$ws.Range("F409").Value = 664339
$ws.Range("G409").Value = 1927

# Add new row 410 for 2021-04-20 (ut 20. 04. 2021)
$ws.Range("A410").Value = 44304
$ws.Range("A410").NumberFormat = $ws.Range("A409").NumberFormat
$ws.Range("B410").Value = 376067
$ws.Range("C410").Value = 1434
$ws.Range("D410").Value = 93
$ws.Range("E410").Value = 11172
$ws.Range("F410").Value = 341058
$ws.Range("G410").Value = 597
